# Applies the cryptos.xlsx data refresh described in the commit:
# "Updated cryptos list on Mon Dec 18 13:48:13 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.350.40"
$ws.Range("E2").Value = "  -1.35%  "

# Row 3
$ws.Range("D3").Value = "2.152.38"
$ws.Range("E3").Value = "  -2.65%  "

# Row 4
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("E5").Value = "  -1.69%  "

# Row 6
$ws.Range("E6").Value = "  -3.58%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "70.32"
$ws.Range("E7").Value = "  -2.47%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("E9").Value = "  -4.72%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.38"
$ws.Range("E10").Value = "  -5.78%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0898"
$ws.Range("E11").Value = "  -5.29%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.11"
$ws.Range("E12").Value = "  -4.65%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0999"
$ws.Range("E13").Value = "  -3.39%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.62"
$ws.Range("E14").Value = "  -4.72%  "

# Row 15
$ws.Range("D15").Value = "2.471.05"
$ws.Range("E15").Value = "  -2.83%  "

# Row 16
$ws.Range("E16").Value = "  -0.05%  "

# Row 17
$ws.Range("D17").Value = "2.144.38"
$ws.Range("E17").Value = "  -2.69%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.782"
$ws.Range("E18").Value = "  -5.77%  "

# Row 19
$ws.Range("D19").Value = "41.176.37"
$ws.Range("E19").Value = "  -1.43%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000100"
$ws.Range("E20").Value = "  -4.35%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.12"
$ws.Range("E21").Value = "  -4.50%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.73"
$ws.Range("E22").Value = "  -6.58%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.64"
$ws.Range("E23").Value = "  -11.68%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "225.29"
$ws.Range("E24").Value = "  -1.63%  "

# Row 25
$ws.Range("E25").Value = "  -3.39%  "

# Row 26
$ws.Range("E26").Value = "  +0.10%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.60"
$ws.Range("E27").Value = "  -7.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.27"
$ws.Range("E28").Value = "  -9.86%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  -3.75%  "

# Row 30
$ws.Range("E30").Value = "  -0.67%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "170.95"

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.66"
$ws.Range("E32").Value = "  -3.59%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.88"
$ws.Range("E33").Value = "  +6.40%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0758"
$ws.Range("E34").Value = "  -4.21%  "

# Row 35
$ws.Range("E35").Value = "  -8.91%  "

# Row 36
$ws.Range("E36").Value = "  -3.44%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.27"
$ws.Range("E37").Value = "  +1.61%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.103"
$ws.Range("E38").Value = "  -3.22%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0290"
$ws.Range("E39").Value = "  -2.75%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.93"
$ws.Range("E40").Value = "  -12.24%  "

# Row 41
$ws.Range("E41").Value = "  -2.34%  "

# Row 42
$ws.Range("E42").Value = "  -6.40%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "58.18"
$ws.Range("E43").Value = "  -9.04%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.188"
$ws.Range("E44").Value = "  -3.97%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.23"
$ws.Range("E45").Value = "  -4.96%  "

# Row 46
$ws.Range("E46").Value = "  -3.96%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.38"
$ws.Range("E47").Value = "  -5.68%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.07"
$ws.Range("E48").Value = "  -2.68%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.10"
$ws.Range("E49").Value = "  -4.50%  "

# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.16"
$ws.Range("E50").Value = "  -6.82%  "

# Row 51
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.61"
$ws.Range("E51").Value = "  -3.16%  "
